# Engineering Offices test data - fill in row 2 (data row) for columns E:J
# (Office Name / Office Number / Secretariat / Municipality / Inspectors / Capacity)
# plus housekeeping: selection cursor, page setup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E2 previously carried the "Hyperlink" look (style index 3); put it back to a
# plain text cell like its neighbours (style index 1) before writing the value.
$ws.Range("E2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"

$ws.Range("E2").Value = "إسم المكتب الهندسى"
$ws.Range("F2").Value = "رقم المكتب الهندسى"
$ws.Range("G2").Value = "الأمانة"
$ws.Range("H2").Value = "البلدية"
$ws.Range("I2").Value = "عدد المراقبين"
$ws.Range("J2").Value = "الطاقة الإستيعابية (عدد الأسّرة)"

# Move/record the active selection like it was left after the edits.
$ws.Range("B13").Select() | Out-Null

# Page setup (printer-independent bits only).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
